$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.775.84"
$ws.Range("E2").Value = "  +0.53%  "
$ws.Range("D3").Value = "1.919.08"
$ws.Range("E3").Value = "  +1.66%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "241.84"
$ws.Range("E5").Value = "  -1.63%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  +0.01%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4930"
$ws.Range("E7").Value = "  +0.41%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3020"
$ws.Range("E8").Value = "  +2.75%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06801"
$ws.Range("E9").Value = "  +0.62%  "
$ws.Range("D10").Value = "1.918.13"
$ws.Range("E10").Value = "  +1.59%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "17.30"
$ws.Range("E11").Value = "  +0.88%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07342"
$ws.Range("E12").Value = "  +1.61%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.222"
$ws.Range("E13").Value = "  +3.36%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "89.34"
$ws.Range("E14").Value = "  -1.83%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6787"
$ws.Range("E15").Value = "  +0.54%  "
$ws.Range("D16").Value = "30.745.09"
$ws.Range("E16").Value = "  +0.55%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008030"
$ws.Range("E17").Value = "  +1.09%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.69"
$ws.Range("E18").Value = "  +4.00%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.001"
$ws.Range("E19").Value = "  +0.01%  "
$ws.Range("D20").Value = "2.164.26"
$ws.Range("E20").Value = "  +1.47%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.399"
$ws.Range("E21").Value = "  +12.23%  "
$ws.Range("E22").Value = "  -0.05%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "202.29"
$ws.Range("E23").Value = "  +11.16%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.355"
$ws.Range("E24").Value = "  +4.91%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.738"
$ws.Range("E25").Value = "  +4.25%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "161.40"
$ws.Range("E26").Value = "  +3.91%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.94"
$ws.Range("E27").Value = "  -0.08%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.975"
$ws.Range("E28").Value = "  +4.05%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.461"
$ws.Range("E29").Value = "  +4.44%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.406"
$ws.Range("E30").Value = "  +2.22%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09202"
$ws.Range("E31").Value = "  +1.94%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.114"
$ws.Range("E32").Value = "  +3.00%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05343"
$ws.Range("E33").Value = "  +2.94%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7527"
$ws.Range("E34").Value = "  -0.02%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.132"
$ws.Range("E35").Value = "  +2.00%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.699"
$ws.Range("E36").Value = "  -1.88%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01874"
$ws.Range("E37").Value = "  +1.43%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.730"
$ws.Range("E38").Value = "  +2.63%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.9325"
$ws.Range("E39").Value = "  -0.35%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.105"
$ws.Range("E40").Value = "  -1.55%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4527"
$ws.Range("E41").Value = "  +2.58%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "73.07"
$ws.Range("E42").Value = "  +26.48%  "
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.990"
$ws.Range("E43").Value = "  +4.51%  "
$ws.Range("B44").Value = "Quant"
$ws.Range("C44").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "107.82"
$ws.Range("E44").Value = "  +2.44%  "
$ws.Range("B45").Value = "Algorand"
$ws.Range("C45").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1406"
$ws.Range("E45").Value = "  +5.28%  "
$ws.Range("B46").Value = "PaxDollar"
$ws.Range("C46").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.003"
$ws.Range("E46").Value = "  +0.20%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.786"
$ws.Range("E47").Value = "  +2.65%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "36.07"
$ws.Range("E48").Value = "  +7.68%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.222"
$ws.Range("E49").Value = "  +6.08%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05964"
$ws.Range("E50").Value = "  +2.04%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4079"
$ws.Range("E51").Value = "  +3.97%  "
